$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 520, shifting all existing rows 520:582 down to 521:583.
$ws.Rows.Item(520).Insert()

# Populate the newly inserted row 520 with the new weekly data record.
$ws.Range("A520").Value = 4
$ws.Range("B520").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C520").Value = "Los Lagos"
$ws.Range("D520").Value = 45154
$ws.Range("E520").Value = 10
$ws.Range("F520").Value = 100112008
$ws.Range("G520").Value = "Coliflor"
$ws.Range("H520").Value = "Sin especificar"
$ws.Range("I520").Value = "Primera"
$ws.Range("J520").Value = 250
$ws.Range("K520").Value = 1500
$ws.Range("L520").Value = 1500
$ws.Range("M520").Value = 1500
$ws.Range("N520").Value = "$/unidad"
$ws.Range("O520").Value = "Región Metropolitana"
$ws.Range("P520").Value = 1500
$ws.Range("Q520").Value = 1
$ws.Range("R520").Value = "Hortaliza"
